$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 36
$templateRow = 35

$values = @(34, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 105, 95, 1426, 1626, 0, 0, 0, 2, 1, 0, 54)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($newRow, $col).Value = $values[$i]
}

# Column A mirrors the bold/bordered/centered style used on the row above
# (column B deliberately keeps the default/general style here, unlike the
# date-formatted B35, since the new value is a plain 0, not a date).
$ws.Range("A$templateRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
